# Swap the values of columns A, B, E, F, G, H, Q, R between rows 2 and 3.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")

    $value2 = $cell2.Value2
    $value3 = $cell3.Value2

    $cell2.Value2 = $value3
    $cell3.Value2 = $value2
}
